$wb = $excel.ActiveWorkbook

# --- TestSteps sheet: mark rows 20-30 (H column) as PASS ---
$steps = $wb.Worksheets.Item("TestSteps")
$steps.Activate()

for ($r = 20; $r -le 30; $r++) {
    $steps.Cells.Item($r, 8).Value = "PASS"
}

# Scroll/selection state on TestSteps: selection now H17:H33, view scrolled down
$steps.Range("A10").Select()
$steps.Range("H17:H33").Select()

# --- TestCases sheet becomes the active tab, selection moves to D4 ---
$cases = $wb.Worksheets.Item("TestCases")
$cases.Activate()
$cases.Range("D4").Select()
